$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. "191.40", "1.00") are not
# auto-converted to numbers by Excel, which would strip trailing zeros.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.760.80'
$ws.Range('E2').Value = '  +9.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.520.40'
$ws.Range('E3').Value = '  +11.14%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '191.40'
$ws.Range('E5').Value = '  +13.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '557.61'
$ws.Range('E6').Value = '  +10.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.520.36'
$ws.Range('E7').Value = '  +11.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').Value = '  +5.63%  '
$ws.Range('E9').Value = '  -0.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.640'
$ws.Range('E10').Value = '  +9.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.11'
$ws.Range('E11').Value = '  +6.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.151'
$ws.Range('E12').Value = '  +18.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000277'
$ws.Range('E13').Value = '  +11.97%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.49'
$ws.Range('E14').Value = '  +9.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.080.55'
$ws.Range('E15').Value = '  +10.67%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.523.88'
$ws.Range('E16').Value = '  +11.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.178.31'
$ws.Range('E17').Value = '  +10.45%  '
$ws.Range('E18').Value = '  +8.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.37'
$ws.Range('E19').Value = '  +10.36%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.87'
$ws.Range('E20').Value = '  +13.02%  '
$ws.Range('E21').Value = '  +9.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '408.33'
$ws.Range('E22').Value = '  +15.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.96'
$ws.Range('E23').Value = '  +10.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.75'
$ws.Range('E24').Value = '  +12.59%  '
$ws.Range('E25').Value = '  +13.29%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '84.85'
$ws.Range('E26').Value = '  +9.28%  '
$ws.Range('E27').Value = '  +14.34%  '
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.92'
$ws.Range('E29').Value = '  +10.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.66'
$ws.Range('E30').Value = '  +9.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.58'
$ws.Range('E31').Value = '  +11.65%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '685.06'
$ws.Range('E32').Value = '  +13.92%  '
$ws.Range('E33').Value = '  +10.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.74'
$ws.Range('E34').Value = '  +9.09%  '
$ws.Range('E35').Value = '  +11.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '60.51'
$ws.Range('E36').Value = '  +8.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0840'
$ws.Range('E37').Value = '  +29.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '39.23'
$ws.Range('E38').Value = '  +11.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.406'
$ws.Range('E39').Value = '  +10.87%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.42'
$ws.Range('E41').Value = '  +28.60%  '
$ws.Range('E42').Value = '  +13.86%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.77'
$ws.Range('E43').Value = '  +19.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.03'
$ws.Range('E44').Value = '  +19.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.064.36'
$ws.Range('E46').Value = '  +10.47%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0422'
$ws.Range('E47').Value = '  +12.92%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.75'
$ws.Range('E48').Value = '  +7.80%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.17'
$ws.Range('E49').Value = '  +9.84%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.95'
$ws.Range('E50').Value = '  +22.80%  '
$ws.Range('E51').Value = '  +9.36%  '
